$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value parses as a plain number need to be forced to remain
# text (matching the source data, which stores these as literal strings) by
# temporarily switching the cell to Text format, then restoring its original style.

$ws.Range("D2").Value = "27.798.21"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.886.11"
$ws.Range("E3").Value = "  +1.56%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  +0.49%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4711"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +2.09%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -0.48%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.65"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +2.44%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08069"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.64%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.026"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +1.44%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.16"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "1.884.47"
$ws.Range("E13").Value = "  +1.84%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.978"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.81%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.126"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.28%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.55%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06745"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +2.46%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.31"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +1.26%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001049"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +1.90%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "27.827.23"
$ws.Range("E22").Value = "  +1.73%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.524"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.75%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.87%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").Value = "2.106.22"
$ws.Range("E26").Value = "  +1.66%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.88"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +3.53%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -0.17%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +1.95%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.581"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +1.96%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.02"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.25%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9827"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("E33").Value = "  +0.72%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.450"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +0.78%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.358"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +1.75%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06163"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +1.99%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02269"
$ws.Range("D38").Style = $origStyle
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.078"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6007"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.40%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1893"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E44").Value = "  -2.00%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5718"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.59%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.25"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +2.02%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.945"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +1.49%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.395"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +2.27%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.17"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +4.60%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000304"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -1.20%  "
